# DiscountReport (row 7) test case: status flipped from "UnderAnalysis" to
# "Passed" and the Manual Status (MnS) column picks up the matching
# green "Passed" formatting used elsewhere in the sheet (e.g. G4/G5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value first ...
$ws.Range("G7").Value = "Passed"

# ... then copy the formatting (fill/font/border) from a cell that already
# carries the "Passed" style so G7 matches G4/G5/B2/B3/B6 visually.
$ws.Range("G4").Copy()
$ws.Range("G7").PasteSpecial(-4122)

# Leave the selection on the cell that was just edited.
$ws.Range("G7").Select()
